$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13
$ws.Range("C2").Value = 15

$ws.Range("B3").Value = 13
$ws.Range("C3").Value = 15

$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 9

$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 13
